$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 55639.277
$ws.Range("I11").Value = 55639.277
$ws.Range("K11").Value = 55639.277
$ws.Range("M11").Value = -55499.277
$ws.Range("H129").Value = 6967.706
$ws.Range("I129").Value = 17275.334
$ws.Range("J129").Value = 1345.3636
$ws.Range("K129").Value = 51826.00199999999
$ws.Range("L129").Value = 4036.0908
$ws.Range("M129").Value = -46826.00199999999
$ws.Range("N129").Value = -14036.0908
$ws.Range("H132").Value = 4314326
$ws.Range("I132").Value = 4468157
$ws.Range("K132").Value = 13404471
$ws.Range("M132").Value = -13401941
$ws.Range("H135").Value = 806.34283
$ws.Range("I135").Value = 669.7931
$ws.Range("K135").Value = 6028.1379
$ws.Range("M135").Value = -3493.1379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 57434.055
$ws.Range("I2").Value = 2071.9092
$ws.Range("J2").Value = 144431.72
$ws.Range("K2").Value = 2071.9092
$ws.Range("L2").Value = 144431.72
$ws.Range("M2").Value = -1958.9092
$ws.Range("N2").Value = -144657.72
$ws.Range("H5").Value = 250117.75
$ws.Range("I5").Value = 333400.34
$ws.Range("J5").Value = 270
$ws.Range("K5").Value = 333400.34
$ws.Range("L5").Value = 270
$ws.Range("M5").Value = -333288.34
$ws.Range("N5").Value = -494
$ws.Range("H61").Value = 1674.7307
$ws.Range("I61").Value = 1512.1052
$ws.Range("J61").Value = 2116.1428
$ws.Range("K61").Value = 1512.1052
$ws.Range("L61").Value = 2116.1428
$ws.Range("M61").Value = -1300.1052
$ws.Range("N61").Value = -2540.1428
$ws.Range("H74").Value = 1388.8485
$ws.Range("I74").Value = 732.6957
$ws.Range("K74").Value = 732.6957
$ws.Range("M74").Value = 141.3043
$ws.Range("H77").Value = 1388.8485
$ws.Range("I77").Value = 732.6957
$ws.Range("K77").Value = 3663.4785
$ws.Range("M77").Value = 704.5214999999998
$ws.Range("H98").Value = 19111
$ws.Range("J98").Value = 19111
$ws.Range("L98").Value = 19111
$ws.Range("N98").Value = -25101
$ws.Range("H102").Value = 61777
$ws.Range("I102").Value = 85904.086
$ws.Range("K102").Value = 85904.086
$ws.Range("M102").Value = -84282.086
$ws.Range("H116").Value = 57434.055
$ws.Range("I116").Value = 2071.9092
$ws.Range("J116").Value = 144431.72
$ws.Range("K116").Value = 2071.9092
$ws.Range("L116").Value = 144431.72
$ws.Range("M116").Value = 222.0907999999999
$ws.Range("N116").Value = -149019.72
$ws.Range("H132").Value = 1905.5264
$ws.Range("I132").Value = 1250.0476
$ws.Range("J132").Value = 2715.2354
$ws.Range("K132").Value = 3750.142800000001
$ws.Range("L132").Value = 8145.706200000001
$ws.Range("M132").Value = -1220.142800000001
$ws.Range("N132").Value = -13205.7062
$ws.Range("H136").Value = 1674.7307
$ws.Range("I136").Value = 1512.1052
$ws.Range("J136").Value = 2116.1428
$ws.Range("K136").Value = 4536.3156
$ws.Range("L136").Value = 6348.428400000001
$ws.Range("M136").Value = -1986.3156
$ws.Range("N136").Value = -11448.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 57434.055
$ws.Range("I3").Value = 2071.9092
$ws.Range("J3").Value = 144431.72
$ws.Range("K3").Value = 2071.9092
$ws.Range("L3").Value = 144431.72
$ws.Range("M3").Value = -1957.9092
$ws.Range("N3").Value = -144659.72
$ws.Range("H4").Value = 250117.75
$ws.Range("I4").Value = 333400.34
$ws.Range("J4").Value = 270
$ws.Range("K4").Value = 333400.34
$ws.Range("L4").Value = 270
$ws.Range("M4").Value = -333285.34
$ws.Range("N4").Value = -500
$ws.Range("H20").Value = 43146.48
$ws.Range("I20").Value = 64413.438
$ws.Range("J20").Value = 5338.5557
$ws.Range("K20").Value = 64413.438
$ws.Range("L20").Value = 5338.5557
$ws.Range("M20").Value = -64166.438
$ws.Range("N20").Value = -5832.5557
$ws.Range("H134").Value = 1997
$ws.Range("I134").Value = 1988
$ws.Range("J134").Value = 2034.5
$ws.Range("K134").Value = 5964
$ws.Range("L134").Value = 6103.5
$ws.Range("M134").Value = -3429
$ws.Range("N134").Value = -11173.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 15000
$ws.Range("J45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -16186
$ws.Range("H105").Value = 2235.625
$ws.Range("I105").Value = 2257.5386
$ws.Range("J105").Value = 2140.6667
$ws.Range("K105").Value = 2257.5386
$ws.Range("L105").Value = 2140.6667
$ws.Range("M105").Value = -510.5385999999999
$ws.Range("N105").Value = -5634.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 111.333336
$ws.Range("I14").Value = 111.333336
$ws.Range("K14").Value = 334.000008
$ws.Range("M14").Value = -161.000008
$ws.Range("H37").Value = 41548.59
$ws.Range("J37").Value = 41548.59
$ws.Range("L37").Value = 124645.77
$ws.Range("N37").Value = -124869.77

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 59142.863
$ws.Range("I70").Value = 99608.664
$ws.Range("J70").Value = 6031.5
$ws.Range("K70").Value = 99608.664
$ws.Range("L70").Value = 6031.5
$ws.Range("M70").Value = -99338.664
$ws.Range("N70").Value = -6571.5
$ws.Range("H73").Value = 59142.863
$ws.Range("I73").Value = 99608.664
$ws.Range("J73").Value = 6031.5
$ws.Range("K73").Value = 99608.664
$ws.Range("L73").Value = 6031.5
$ws.Range("M73").Value = -98672.664
$ws.Range("N73").Value = -7903.5
$ws.Range("H122").Value = 944.7778
$ws.Range("I122").Value = 953.73334
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 2861.20002
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -411.2000200000002
$ws.Range("N122").Value = -7600
$ws.Range("H132").Value = 1423
$ws.Range("I132").Value = 1040.125
$ws.Range("J132").Value = 2129.8462
$ws.Range("K132").Value = 3120.375
$ws.Range("L132").Value = 6389.5386
$ws.Range("M132").Value = -590.375
$ws.Range("N132").Value = -11449.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2992.9412
$ws.Range("I7").Value = 1934.5454
$ws.Range("J7").Value = 4933.3335
$ws.Range("K7").Value = 1934.5454
$ws.Range("L7").Value = 4933.3335
$ws.Range("M7").Value = -1822.5454
$ws.Range("N7").Value = -5157.3335
$ws.Range("H40").Value = 60735.176
$ws.Range("I40").Value = 200979.8
$ws.Range("J40").Value = 2299.9167
$ws.Range("K40").Value = 200979.8
$ws.Range("L40").Value = 2299.9167
$ws.Range("M40").Value = -200843.8
$ws.Range("N40").Value = -2571.9167
$ws.Range("H46").Value = 1446636
$ws.Range("I46").Value = 395
$ws.Range("J46").Value = 2025132.4
$ws.Range("K46").Value = 395
$ws.Range("L46").Value = 2025132.4
$ws.Range("M46").Value = -207
$ws.Range("N46").Value = -2025508.4
$ws.Range("H93").Value = 1425.8334
$ws.Range("I93").Value = 1391.25
$ws.Range("J93").Value = 1495
$ws.Range("K93").Value = 1391.25
$ws.Range("L93").Value = 1495
$ws.Range("M93").Value = -143.25
$ws.Range("N93").Value = -3991
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 2992.9412
$ws.Range("I126").Value = 1934.5454
$ws.Range("J126").Value = 4933.3335
$ws.Range("K126").Value = 5803.6362
$ws.Range("L126").Value = 14800.0005
$ws.Range("M126").Value = -3333.6362
$ws.Range("N126").Value = -19740.0005
$ws.Range("H132").Value = 2849.2368
$ws.Range("I132").Value = 3366.8262
$ws.Range("J132").Value = 2055.6
$ws.Range("K132").Value = 10100.4786
$ws.Range("L132").Value = 6166.799999999999
$ws.Range("M132").Value = -7570.4786
$ws.Range("N132").Value = -11226.8
$ws.Range("H136").Value = 1305.4103
$ws.Range("I136").Value = 1235.0303
$ws.Range("J136").Value = 1692.5
$ws.Range("K136").Value = 3705.0909
$ws.Range("L136").Value = 5077.5
$ws.Range("M136").Value = -1155.0909
$ws.Range("N136").Value = -10177.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 36215.8
$ws.Range("I56").Value = 8900
$ws.Range("J56").Value = 39250.89
$ws.Range("K56").Value = 8900
$ws.Range("L56").Value = 39250.89
$ws.Range("M56").Value = -8186
$ws.Range("N56").Value = -40678.89
$ws.Range("H101").Value = 12020.4
$ws.Range("J101").Value = 12020.4
$ws.Range("L101").Value = 12020.4
$ws.Range("N101").Value = -18510.4
$ws.Range("H122").Value = 1749.5
$ws.Range("I122").Value = 998.625
$ws.Range("K122").Value = 2995.875
$ws.Range("M122").Value = -545.875
$ws.Range("H126").Value = 1685.2667
$ws.Range("I126").Value = 1643.8182
$ws.Range("J126").Value = 1799.25
$ws.Range("K126").Value = 4931.4546
$ws.Range("L126").Value = 5397.75
$ws.Range("M126").Value = -2461.4546
$ws.Range("N126").Value = -10337.75
$ws.Range("H132").Value = 1904.2858
$ws.Range("J132").Value = 3578.6924
$ws.Range("L132").Value = 10736.0772
$ws.Range("N132").Value = -15796.0772
$ws.Range("H136").Value = 470.70212
$ws.Range("I136").Value = 288.76315
$ws.Range("J136").Value = 1238.8889
$ws.Range("K136").Value = 866.28945
$ws.Range("L136").Value = 3716.6667
$ws.Range("M136").Value = 1683.71055
$ws.Range("N136").Value = -8816.6667
